$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Last existing data row is 357 (date serial 44431 = 2021-08-23).
# Append rows 358..366 continuing the daily date series through
# 44440 (2021-09-01), matching the commit "aggiornamento fino a 1/09/2021".
$lastRow = 357
$lastDate = 44431
$newRowsCount = 9

for ($i = 1; $i -le $newRowsCount; $i++) {
    $row = $lastRow + $i
    $dateValue = $lastDate + $i

    # Replicate column A's existing date style (borders, bold, centered,
    # custom date number format) from the previous row instead of building
    # a brand-new style.
    $ws.Range("A$lastRow").Copy()
    $ws.Range("A$row").PasteSpecial(-4122)

    $ws.Range("A$row").Value = $dateValue
    $ws.Range("B$row").Value = 0
    $ws.Range("C$row").Value = 0
    $ws.Range("D$row").Value = 0
}
